# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values change for rows 2-15 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 3
    4  = 4
    5  = 3
    6  = 6
    7  = 5
    8  = 2
    9  = 5
    10 = 5
    11 = 2
    12 = 2
    13 = 1
    14 = 1
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
